$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 corresponds to Conor McGregor - add train/val set paths in E12/F12
$ws.Range("E12").Value = "data/european-celebrities/train/connor_mcgregor"
$ws.Range("F12").Value = "data/european-celebrities/val/connor_mcgregor"

# Row 5 corresponds to Adele - add train/val set paths in E5/F5
$ws.Range("E5").Value = "data/european-celebrities/train/adele"
$ws.Range("F5").Value = "data/european-celebrities/val/adele"

# Update the active cell selection to F6 as in the edited workbook
$ws.Range("F6").Select()
